$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-17 Tuesday" "2026-02-18 Wednesday"

Replace-Text "827×6=4962" "695×4=2780"
Replace-Text "934×9=8406" "935×5=4675"
Replace-Text "129×3=387" "745×8=5960"
Replace-Text "265×4=1060" "674×8=5392"
Replace-Text "739×8=5912" "437×2=874"
Replace-Text "866×8=6928" "602×2=1204"
Replace-Text "960×9=8640" "297×5=1485"
Replace-Text "744×6=4464" "369×7=2583"
Replace-Text "975×8=7800" "862×9=7758"
Replace-Text "225×3=675" "710×5=3550"
Replace-Text "678×9=6102" "341×5=1705"
Replace-Text "527×9=4743" "911×9=8199"
Replace-Text "913×6=5478" "562×6=3372"
Replace-Text "898×6=5388" "770×6=4620"
Replace-Text "867×5=4335" "651×4=2604"
Replace-Text "893×4=3572" "271×5=1355"
Replace-Text "987×5=4935" "455×2=910"
Replace-Text "208×8=1664" "608×5=3040"
Replace-Text "552×4=2208" "163×9=1467"
Replace-Text "171×7=1197" "308×4=1232"
Replace-Text "781×2=1562" "452×2=904"
Replace-Text "447×4=1788" "799×6=4794"
Replace-Text "900×8=7200" "838×4=3352"
Replace-Text "538×2=1076" "186×2=372"
Replace-Text "795×2=1590" "299×2=598"
